$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains text formatting so numeric-looking price strings
# (e.g. "0.9994", "334.46") are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.124.13'
$ws.Range('E2').Value = '  +2.63%  '
$ws.Range('D3').Value = '1.917.96'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '334.46'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = '0.9991'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '0.4668'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '0.4094'
$ws.Range('E8').Value = '  +3.54%  '
$ws.Range('D9').Value = '48.26'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').Value = '0.08004'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').Value = '1.011'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').Value = '21.93'
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '1.910.29'
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('D14').Value = '5.983'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '7.176'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '89.93'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '0.9999'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.00001034'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.06581'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('D20').Value = '17.63'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '0.9980'
$ws.Range('E21').Value = '  -0.36%  '
$ws.Range('D22').Value = '29.116.29'
$ws.Range('E22').Value = '  +2.57%  '
$ws.Range('D23').Value = '5.493'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D24').Value = '11.16'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = '2.239'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').Value = '2.137.97'
$ws.Range('E26').Value = '  +2.65%  '
$ws.Range('D27').Value = '157.49'
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('D28').Value = '19.83'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').Value = '2.122'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '5.415'
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').Value = '119.80'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('D32').Value = '0.9880'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('D33').Value = '0.09429'
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.430'
$ws.Range('E34').Value = '  +3.61%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '3.599'
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('D36').Value = '5.351'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').Value = '0.06120'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '0.02239'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').Value = '8.448'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').Value = '1.171'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.5876'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '0.9977'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '10.25'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.1835'
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('D45').Value = '1.258'
$ws.Range('E45').Value = '  -3.16%  '
$ws.Range('D46').Value = '2.350'
$ws.Range('E46').Value = '  +16.63%  '
$ws.Range('D47').Value = '12.15'
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('D48').Value = '0.5539'
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.930'
$ws.Range('E49').Value = '  -1.61%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.07109'
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '47.50'
$ws.Range('E51').Value = '  +23.21%  '
